$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.015397587163892
$ws.Range("E2").Value = -0.001911419365408951
$ws.Range("D3").Value = 0.05235367123662581
$ws.Range("E3").Value = -0.01575838335286517
$ws.Range("D4").Value = 0.01496197357468868
$ws.Range("E4").Value = 0.007909910284964683
$ws.Range("D5").Value = 0.00976183051404096
$ws.Range("E5").Value = -0.003190403266972952
$ws.Range("D6").Value = 0.01556374134265798
$ws.Range("E6").Value = -0.006191318546740932
$ws.Range("D7").Value = 0.02003177825892858
$ws.Range("E7").Value = -0.01622214103032515
$ws.Range("D8").Value = 0.003836744469487775
$ws.Range("E8").Value = -0.01005738962441061
$ws.Range("D9").Value = 0.006278665874268305
$ws.Range("E9").Value = -0.01303571428571426
$ws.Range("D10").Value = 0.01408719596022216
$ws.Range("E10").Value = -0.008794619762027889
$ws.Range("D11").Value = 0.00881193314715017
$ws.Range("E11").Value = -0.008058245564430577
$ws.Range("D12").Value = 0.01396255696801571
$ws.Range("E12").Value = -0.0133832976445396
$ws.Range("D13").Value = 0.002979096151801263
$ws.Range("E13").Value = -0.01536772777167938
$ws.Range("D14").Value = 0.005943865982659257
$ws.Range("E14").Value = -0.01498559077809802
$ws.Range("D15").Value = 0.01359648832373917
$ws.Range("E15").Value = -0.02105752623887325
$ws.Range("D16").Value = 0.009901512143237853
$ws.Range("E16").Value = -0.01188959660297251
$ws.Range("D17").Value = 0.02166279875414402
$ws.Range("E17").Value = 0.0005822604645144835
$ws.Range("D18").Value = 0.008675023349614041
$ws.Range("E18").Value = 0
$ws.Range("D19").Value = 0.0166025864761009
$ws.Range("E19").Value = -0.008463893390959876
$ws.Range("D20").Value = 0.01189396080051124
$ws.Range("E20").Value = -0.0204241948153967
$ws.Range("D21").Value = 0.006939358423160804
$ws.Range("E21").Value = 0.04151444702756568
$ws.Range("D22").Value = 0.01337704645989942
$ws.Range("E22").Value = -0.006779661016949157
$ws.Range("D23").Value = 0.01907312937859201
$ws.Range("E23").Value = -0.003659289304293933
$ws.Range("D24").Value = 0.009765194085045033
$ws.Range("E24").Value = -0.05660377358490565
$ws.Range("D25").Value = 0.02050769241182885
$ws.Range("E25").Value = -0.003158809128958384
$ws.Range("D26").Value = 0.01167897255383428
$ws.Range("E26").Value = -0.005936047488379836
$ws.Range("D27").Value = 0.02196832312034725
$ws.Range("E27").Value = -0.04165302144249516
$ws.Range("D28").Value = 0.05820846487602908
$ws.Range("E28").Value = -0.01168539325842699
$ws.Range("D29").Value = 0.02175031388860182
$ws.Range("E29").Value = -0.01427027027027028
$ws.Range("D30").Value = 0.03103330269720194
$ws.Range("E30").Value = -0.02211874272409797
$ws.Range("D31").Value = 0.01569364370060229
$ws.Range("E31").Value = -0.02364343931272472
$ws.Range("D32").Value = 0.01381593018757893
$ws.Range("E32").Value = -0.00702415624464614
$ws.Range("D33").Value = 0.02006706461029538
$ws.Range("E33").Value = -0.05341378925019558
$ws.Range("D34").Value = 0.04257439998404521
$ws.Range("E34").Value = -0.01133715188623341
$ws.Range("D35").Value = 0.01101008908666335
$ws.Range("E35").Value = -0.01731160896130357
$ws.Range("D36").Value = 0.009662169150957337
$ws.Range("E36").Value = -0.006768953068592043
$ws.Range("D37").Value = 0.01067466631153503
$ws.Range("E37").Value = -0.02297592997811815
$ws.Range("D38").Value = 0.007435360569557616
$ws.Range("E38").Value = -0.006785624528775958
$ws.Range("D39").Value = 0.0112148931878002
$ws.Range("E39").Value = -0.009297520661157077
$ws.Range("D40").Value = 0.0171435608125892
$ws.Range("E40").Value = -0.008280377431157282
$ws.Range("D41").Value = 0.01674740689433179
$ws.Range("E41").Value = -0.0007773302240122737
$ws.Range("D42").Value = 0.03398047606864048
$ws.Range("E42").Value = -0.001484780994803314
$ws.Range("D43").Value = 0.01123865681694397
$ws.Range("E43").Value = -0.006416189901530034
$ws.Range("D44").Value = 0.02195960275107743
$ws.Range("E44").Value = 0.0005275884809847753
$ws.Range("D45").Value = 0.01288780259968672
$ws.Range("E45").Value = -0.01017855977728976
$ws.Range("D46").Value = 0.007799560563276333
$ws.Range("E46").Value = -0.006181269316466675
$ws.Range("D47").Value = 0.0130315952601108
$ws.Range("E47").Value = -0.01535748085692168
$ws.Range("D48").Value = 0.009615608607891705
$ws.Range("E48").Value = -0.03130738959335366
$ws.Range("D49").Value = 0.01516465987185999
$ws.Range("E49").Value = -0.01997461584907523
$ws.Range("D50").Value = 0.008788357005945701
$ws.Range("E50").Value = -0.008795710584974858
$ws.Range("D51").Value = 0.01104185614614626
$ws.Range("E51").Value = 0.01736334405144691
$ws.Range("D52").Value = 0.008843513341577293
$ws.Range("E52").Value = -0.003912605562872873
$ws.Range("D53").Value = 0.009689420304925514
$ws.Range("E53").Value = 0.0004724940938238831
$ws.Range("D54").Value = 0.1352698698065873
$ws.Range("E54").Value = 0.0001970831690973895
$ws.Range("D55").Value = 0.04404502797304794
$ws.Range("E55").Value = -0.007092198581560405
$ws.Range("E56").Value = -0.009881123297001526

$ws.Protect()
